# Auto-generated edit script applying the Golem_Profits scheduled-runner update.
# For each touched leve row, recompute currentAveragePrice / NQ / HQ price+profit columns.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1501.2727
$ws.Range("I28").Value = 635
$ws.Range("J28").Value = 2223.1667
$ws.Range("K28").Value = 635
$ws.Range("L28").Value = 2223.1667
$ws.Range("M28").Value = -150
$ws.Range("N28").Value = -3193.1667
$ws.Range("H39").Value = 462.08334
$ws.Range("I39").Value = 549.7778
$ws.Range("K39").Value = 1649.3334
$ws.Range("M39").Value = -1353.3334
$ws.Range("H76").Value = 3949.5
$ws.Range("I76").Value = 3949.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3949.5
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3634.5
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3949.5
$ws.Range("I79").Value = 3949.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3949.5
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2857.5
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 667.3333
$ws.Range("I98").Value = 291
$ws.Range("J98").Value = 1420
$ws.Range("K98").Value = 291
$ws.Range("L98").Value = 1420
$ws.Range("M98").Value = 1207
$ws.Range("N98").Value = -4416
$ws.Range("H100").Value = 1912.75
$ws.Range("I100").Value = 1216.3334
$ws.Range("K100").Value = 1216.3334
$ws.Range("M100").Value = -675.3334
$ws.Range("H122").Value = 667.3333
$ws.Range("I122").Value = 291
$ws.Range("J122").Value = 1420
$ws.Range("K122").Value = 873
$ws.Range("L122").Value = 4260
$ws.Range("M122").Value = 1577
$ws.Range("N122").Value = -9160

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -885
$ws.Range("H30").Value = 7888.25
$ws.Range("I30").Value = 7934.3335
$ws.Range("J30").Value = 7750
$ws.Range("K30").Value = 7934.3335
$ws.Range("L30").Value = 7750
$ws.Range("M30").Value = -7784.3335
$ws.Range("N30").Value = -8050
$ws.Range("H97").Value = 66669350
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H122").Value = 1414.8334
$ws.Range("I122").Value = 1172.25
$ws.Range("K122").Value = 3516.75
$ws.Range("M122").Value = -1066.75
$ws.Range("H124").Value = 72356.75
$ws.Range("J124").Value = 72356.75
$ws.Range("L124").Value = 72356.75
$ws.Range("N124").Value = -82176.75
$ws.Range("H135").Value = 19999.5
$ws.Range("J135").Value = 19999.5
$ws.Range("L135").Value = 19999.5
$ws.Range("N135").Value = -30139.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H94").Value = 158565.58
$ws.Range("I94").Value = 158565.58
$ws.Range("K94").Value = 158565.58
$ws.Range("M94").Value = -158114.58
$ws.Range("H99").Value = 4894.7144
$ws.Range("I99").Value = 4925.3335
$ws.Range("J99").Value = 4711
$ws.Range("K99").Value = 4925.3335
$ws.Range("L99").Value = 4711
$ws.Range("M99").Value = -3427.3335
$ws.Range("N99").Value = -7707
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 845.2
$ws.Range("I35").Value = 799.5
$ws.Range("K35").Value = 799.5
$ws.Range("M35").Value = -505.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1810.7778
$ws.Range("J68").Value = 1505.8
$ws.Range("L68").Value = 4517.4
$ws.Range("N68").Value = -6139.4
$ws.Range("H71").Value = 1810.7778
$ws.Range("J71").Value = 1505.8
$ws.Range("L71").Value = 13552.2
$ws.Range("N71").Value = -21664.2
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 9000
$ws.Range("M76").Value = -8617
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 9000
$ws.Range("M79").Value = -7674
$ws.Range("H131").Value = 1401.7
$ws.Range("J131").Value = 3100
$ws.Range("L131").Value = 9300
$ws.Range("N131").Value = -19380
$ws.Range("H141").Value = 9996.666999999999
$ws.Range("I141").Value = 9997.5
$ws.Range("K141").Value = 29992.5
$ws.Range("M141").Value = -24812.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 503.2
$ws.Range("I113").Value = 503.2
$ws.Range("K113").Value = 503.2
$ws.Range("M113").Value = 1666.8
$ws.Range("H122").Value = 3972.6667
$ws.Range("I122").Value = 3548.7368
$ws.Range("K122").Value = 10646.2104
$ws.Range("M122").Value = -8196.2104

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1660.3334
$ws.Range("I32").Value = 1660.3334
$ws.Range("K32").Value = 1660.3334
$ws.Range("M32").Value = -1343.3334
$ws.Range("H64").Value = 37459.8
$ws.Range("J64").Value = 37459.8
$ws.Range("L64").Value = 37459.8
$ws.Range("N64").Value = -37909.8
$ws.Range("H67").Value = 37459.8
$ws.Range("J67").Value = 37459.8
$ws.Range("L67").Value = 37459.8
$ws.Range("N67").Value = -39019.8
$ws.Range("H70").Value = 30163
$ws.Range("J70").Value = 30163
$ws.Range("L70").Value = 30163
$ws.Range("N70").Value = -30703
$ws.Range("H73").Value = 30163
$ws.Range("J73").Value = 30163
$ws.Range("L73").Value = 30163
$ws.Range("N73").Value = -32035
$ws.Range("H93").Value = 37037760
$ws.Range("I93").Value = 41667350
$ws.Range("K93").Value = 41667350
$ws.Range("M93").Value = -41666102
$ws.Range("H100").Value = 2135.375
$ws.Range("I100").Value = 2135.375
$ws.Range("K100").Value = 2135.375
$ws.Range("M100").Value = -1594.375
$ws.Range("H122").Value = 3634.5862
$ws.Range("I122").Value = 3093.0715
$ws.Range("K122").Value = 9279.2145
$ws.Range("M122").Value = -6829.2145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376
$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880
$ws.Range("H96").Value = 1646.1111
$ws.Range("I96").Value = 1183.2
$ws.Range("J96").Value = 2224.75
$ws.Range("K96").Value = 1183.2
$ws.Range("L96").Value = 2224.75
$ws.Range("M96").Value = 189.8
$ws.Range("N96").Value = -4970.75
$ws.Range("H132").Value = 2748.75
